$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextCell "D2" "65.197.20"
Set-TextCell "E2" "  +1.94%  "

Set-TextCell "D3" "3.172.33"
Set-TextCell "E3" "  +3.91%  "

Set-TextCell "E4" "  +0.09%  "

Set-TextCell "D5" "573.67"

Set-TextCell "D6" "151.09"
Set-TextCell "E6" "  +6.05%  "

Set-TextCell "E7" "  +0.01%  "

Set-TextCell "D8" "3.170.87"
Set-TextCell "E8" "  +3.99%  "

Set-TextCell "D9" "0.527"
Set-TextCell "E9" "  +3.06%  "

Set-TextCell "E10" "  +5.34%  "

Set-TextCell "E11" "  +2.16%  "

Set-TextCell "E13" "  +18.87%  "

Set-TextCell "D14" "38.21"
Set-TextCell "E14" "  +8.91%  "

Set-TextCell "D15" "3.690.63"
Set-TextCell "E15" "  +4.05%  "

Set-TextCell "D16" "65.299.08"
Set-TextCell "E16" "  +2.04%  "

Set-TextCell "B17" "WrappedEther"
Set-TextCell "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D17" "3.183.39"
Set-TextCell "E17" "  +4.26%  "

Set-TextCell "B18" "Polkadot"
Set-TextCell "C18" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D18" "7.22"
Set-TextCell "E18" "  +7.27%  "

Set-TextCell "E19" "  +1.39%  "

Set-TextCell "D20" "511.94"
Set-TextCell "E20" "  +7.31%  "

Set-TextCell "D21" "14.97"
Set-TextCell "E21" "  +7.01%  "

Set-TextCell "D22" "16.10"
Set-TextCell "E22" "  +13.27%  "

Set-TextCell "D23" "0.737"
Set-TextCell "E23" "  +8.74%  "

Set-TextCell "E24" "  +3.92%  "

Set-TextCell "D25" "84.99"
Set-TextCell "E25" "  +3.85%  "

Set-TextCell "E26" "  -0.07%  "

Set-TextCell "D27" "9.14"
Set-TextCell "E27" "  +15.38%  "

Set-TextCell "E28" "  +4.42%  "

Set-TextCell "E29" "  +9.27%  "

Set-TextCell "D30" "28.06"
Set-TextCell "E30" "  +6.97%  "

Set-TextCell "D31" "2.80"
Set-TextCell "E31" "  +15.25%  "

Set-TextCell "E32" "  +7.90%  "

Set-TextCell "E33" "  +0.09%  "

Set-TextCell "E34" "  +12.45%  "

Set-TextCell "D35" "6.69"
Set-TextCell "E35" "  +8.07%  "

Set-TextCell "D36" "55.61"
Set-TextCell "E36" "  +1.63%  "

Set-TextCell "D37" "479.38"
Set-TextCell "E37" "  +8.41%  "

Set-TextCell "B38" "Hedera"
Set-TextCell "C38" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D38" "0.0881"
Set-TextCell "E38" "  +9.31%  "

Set-TextCell "B39" "dogwifhat"
Set-TextCell "C39" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D39" "3.13"
Set-TextCell "E39" "  +10.50%  "

Set-TextCell "E40" "  +4.06%  "

Set-TextCell "D41" "3.135.73"
Set-TextCell "E41" "  +5.49%  "

Set-TextCell "D42" "8.65"
Set-TextCell "E42" "  +5.16%  "

Set-TextCell "E43" "  +7.63%  "

Set-TextCell "D44" "2.52"
Set-TextCell "E44" "  +17.06%  "

Set-TextCell "E45" "  +11.62%  "

Set-TextCell "D46" "29.23"
Set-TextCell "E46" "  +5.70%  "

Set-TextCell "D47" "0.0₃0587"
Set-TextCell "E47" "  +14.19%  "

Set-TextCell "E48" "  -0.06%  "

Set-TextCell "E49" "  +2.40%  "

Set-TextCell "D50" "2.32"
Set-TextCell "E50" "  +12.62%  "

Set-TextCell "D51" "123.51"
Set-TextCell "E51" "  +5.47%  "

